$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (row 6); dimension will shrink to A1:AH5 automatically
$ws.Rows(6).Delete()

# Update data rows 2-5 with the refreshed dataset / custom-accuracy values
# Row 2
$ws.Cells.Item(2, 1).Value = 45038.50694444445
$ws.Cells.Item(2, 2).Value = 5.237
$ws.Cells.Item(2, 3).Value = 1.607
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 2.478
$ws.Cells.Item(2, 6).Value = 3.053
$ws.Cells.Item(2, 7).Value = 2.519
$ws.Cells.Item(2, 8).Value = 5.488
$ws.Cells.Item(2, 9).Value = 1.554
$ws.Cells.Item(2, 10).Value = 0.9409999999999999
$ws.Cells.Item(2, 11).Value = 4.022
$ws.Cells.Item(2, 12).Value = 1.069
$ws.Cells.Item(2, 13).Value = 0.9379999999999999
$ws.Cells.Item(2, 14).Value = 0.6929999999999999
$ws.Cells.Item(2, 15).Value = 0.87
$ws.Cells.Item(2, 16).Value = 2.749
$ws.Cells.Item(2, 17).Value = 1.106
$ws.Cells.Item(2, 18).Value = 0.51
$ws.Cells.Item(2, 19).Value = 0.063
$ws.Cells.Item(2, 20).Value = 20.094
$ws.Cells.Item(2, 21).Value = 4.803
$ws.Cells.Item(2, 22).Value = 2.372
$ws.Cells.Item(2, 23).Value = 3.891
$ws.Cells.Item(2, 24).Value = 1.034
$ws.Cells.Item(2, 25).Value = 0.249
$ws.Cells.Item(2, 26).Value = 1.98
$ws.Cells.Item(2, 27).Value = 1.154
$ws.Cells.Item(2, 28).Value = 0.674
$ws.Cells.Item(2, 29).Value = 0.9399999999999999
$ws.Cells.Item(2, 30).Value = 3.129
$ws.Cells.Item(2, 31).Value = 2.866
$ws.Cells.Item(2, 32).Value = 3.128
$ws.Cells.Item(2, 33).Value = 0.446
$ws.Cells.Item(2, 34).Value = 1.43

# Row 3
$ws.Cells.Item(3, 1).Value = 45038.51388888889
$ws.Cells.Item(3, 2).Value = 2.132
$ws.Cells.Item(3, 3).Value = 0.8169999999999999
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = 0.712
$ws.Cells.Item(3, 6).Value = 1.299
$ws.Cells.Item(3, 7).Value = 1.235
$ws.Cells.Item(3, 8).Value = 3.715
$ws.Cells.Item(3, 9).Value = 0.467
$ws.Cells.Item(3, 10).Value = 0.5
$ws.Cells.Item(3, 11).Value = 1.665
$ws.Cells.Item(3, 12).Value = 0.429
$ws.Cells.Item(3, 13).Value = 0.18
$ws.Cells.Item(3, 14).Value = 0.269
$ws.Cells.Item(3, 15).Value = 0.246
$ws.Cells.Item(3, 16).Value = 1.507
$ws.Cells.Item(3, 17).Value = 0.361
$ws.Cells.Item(3, 18).Value = 0.348
$ws.Cells.Item(3, 19).Value = 0
$ws.Cells.Item(3, 20).Value = 3.735
$ws.Cells.Item(3, 21).Value = 2.189
$ws.Cells.Item(3, 22).Value = 1.028
$ws.Cells.Item(3, 23).Value = 1.749
$ws.Cells.Item(3, 24).Value = 0.441
$ws.Cells.Item(3, 25).Value = 0.099
$ws.Cells.Item(3, 26).Value = 1.782
$ws.Cells.Item(3, 27).Value = 0.522
$ws.Cells.Item(3, 28).Value = 0.156
$ws.Cells.Item(3, 29).Value = 0.288
$ws.Cells.Item(3, 30).Value = 1.482
$ws.Cells.Item(3, 31).Value = 1.133
$ws.Cells.Item(3, 32).Value = 2.803
$ws.Cells.Item(3, 33).Value = 0.099
$ws.Cells.Item(3, 34).Value = 0.525

# Row 4
$ws.Cells.Item(4, 1).Value = 45038.52083333334
$ws.Cells.Item(4, 2).Value = 17.072
$ws.Cells.Item(4, 3).Value = 12.381
$ws.Cells.Item(4, 4).Value = 0.472
$ws.Cells.Item(4, 5).Value = 34.63
$ws.Cells.Item(4, 6).Value = 29.123
$ws.Cells.Item(4, 7).Value = 13.209
$ws.Cells.Item(4, 8).Value = 42.94
$ws.Cells.Item(4, 9).Value = 19.289
$ws.Cells.Item(4, 10).Value = 8.832000000000001
$ws.Cells.Item(4, 11).Value = 13.823
$ws.Cells.Item(4, 12).Value = 14.03
$ws.Cells.Item(4, 13).Value = 14.589
$ws.Cells.Item(4, 14).Value = 4.108
$ws.Cells.Item(4, 15).Value = 12.453
$ws.Cells.Item(4, 16).Value = 18.556
$ws.Cells.Item(4, 17).Value = 10.401
$ws.Cells.Item(4, 18).Value = 0.403
$ws.Cells.Item(4, 19).Value = 0.428
$ws.Cells.Item(4, 20).Value = 186.957
$ws.Cells.Item(4, 21).Value = 35.535
$ws.Cells.Item(4, 22).Value = 12.022
$ws.Cells.Item(4, 23).Value = 24.166
$ws.Cells.Item(4, 24).Value = 12.456
$ws.Cells.Item(4, 25).Value = 1.676
$ws.Cells.Item(4, 26).Value = 21.977
$ws.Cells.Item(4, 27).Value = 10.382
$ws.Cells.Item(4, 28).Value = 8.887
$ws.Cells.Item(4, 29).Value = 10.535
$ws.Cells.Item(4, 30).Value = 15.41
$ws.Cells.Item(4, 31).Value = 0.716
$ws.Cells.Item(4, 32).Value = 38.297
$ws.Cells.Item(4, 33).Value = 6.525
$ws.Cells.Item(4, 34).Value = 14.522

# Row 5
$ws.Cells.Item(5, 1).Value = 45038.52777777778
$ws.Cells.Item(5, 2).Value = 9.67
$ws.Cells.Item(5, 3).Value = 6.96
$ws.Cells.Item(5, 4).Value = 0.23
$ws.Cells.Item(5, 5).Value = 19.18
$ws.Cells.Item(5, 6).Value = 16.23
$ws.Cells.Item(5, 7).Value = 7.46
$ws.Cells.Item(5, 8).Value = 30.17
$ws.Cells.Item(5, 9).Value = 10.67
$ws.Cells.Item(5, 10).Value = 4.97
$ws.Cells.Item(5, 11).Value = 7.81
$ws.Cells.Item(5, 12).Value = 7.82
$ws.Cells.Item(5, 13).Value = 8.039999999999999
$ws.Cells.Item(5, 14).Value = 2.3
$ws.Cells.Item(5, 15).Value = 6.9
$ws.Cells.Item(5, 16).Value = 10.5
$ws.Cells.Item(5, 17).Value = 5.77
$ws.Cells.Item(5, 18).Value = 0.29
$ws.Cells.Item(5, 19).Value = 0.21
$ws.Cells.Item(5, 20).Value = 101.23
$ws.Cells.Item(5, 21).Value = 20.08
$ws.Cells.Item(5, 22).Value = 6.76
$ws.Cells.Item(5, 23).Value = 13.66
$ws.Cells.Item(5, 24).Value = 6.94
$ws.Cells.Item(5, 25).Value = 0.9399999999999999
$ws.Cells.Item(5, 26).Value = 14.6
$ws.Cells.Item(5, 27).Value = 5.8
$ws.Cells.Item(5, 28).Value = 4.92
$ws.Cells.Item(5, 29).Value = 5.85
$ws.Cells.Item(5, 30).Value = 8.699999999999999
$ws.Cells.Item(5, 31).Value = 0.52
$ws.Cells.Item(5, 32).Value = 27.38
$ws.Cells.Item(5, 33).Value = 3.6
$ws.Cells.Item(5, 34).Value = 8.07

# Column width tweaks (stored width = ColumnWidth + 5/6)
$ws.Columns(5).ColumnWidth = 6.166666666666667
$ws.Columns(7).ColumnWidth = 7.166666666666667
$ws.Columns(8).ColumnWidth = 6.166666666666667
$ws.Columns(12).ColumnWidth = 6.166666666666667
$ws.Columns(27).ColumnWidth = 7.166666666666667
$ws.Columns(30).ColumnWidth = 6.166666666666667
